$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows appended to the TATAMOTORS.NS basket table (rows 613-633).
# Column A holds a literal "YYYY-MM-DD" text date (like the rest of the sheet),
# so the target range is pre-formatted as Text to stop Excel from silently
# re-interpreting the string as a date serial number.
$ws.Range("A613:A633").NumberFormat = "@"

$newRows = @(
    @{ Row=613; Date="2024-08-28"; C=1074.550048828125; D=700.0999755859375; E=260.4500122070312; F=491.7000122070312; G=1535.849975585938; H=29293.15045166016; I=0; J=261.8665328852292 },
    @{ Row=614; Date="2024-08-29"; C=1121.650024414062; D=689.1500244140625; E=253.8000030517578; F=490.2000122070312; G=1520.25; H=29353.35063171387; I=0.002055094079179154; J=262.4046932464968 },
    @{ Row=615; Date="2024-08-30"; C=1111.349975585938; D=683.3499755859375; E=256.4500122070312; F=492.8999938964844; G=1521; H=29291.19976806641; I=-0.002117334556700047; J=261.8490947216458 },
    @{ Row=616; Date="2024-09-02"; C=1092.650024414062; D=715.0499877929688; E=251.3500061035156; F=490.5; G=1505.25; H=29322.30020141602; I=0.00106176713811209; J=262.1271174855656 },
    @{ Row=617; Date="2024-09-03"; C=1085.099975585938; D=710.7999877929688; E=251; F=488.8500061035156; G=1509; H=29201.39971923828; I=-0.004123158188384413; J=261.0463259147074 },
    @{ Row=618; Date="2024-09-04"; C=1080.449951171875; D=722.4000244140625; E=250.5; F=484.1499938964844; G=1488.099975585938; H=29200.89978027344; I=-0.00001712037675078922; J=261.0418567032584 },
    @{ Row=619; Date="2024-09-05"; C=1069.150024414062; D=733.8499755859375; E=251.1499938964844; F=495.6499938964844; G=1447.599975585938; H=29262.3996887207; I=0.002106096350113556; J=261.591636004888 },
    @{ Row=620; Date="2024-09-06"; C=1049.349975585938; D=718.9000244140625; E=247.8000030517578; F=483; G=1418.050048828125; H=28702.20024108887; I=-0.01914400232349252; J=256.5837251174042 },
    @{ Row=621; Date="2024-09-09"; C=1038.699951171875; D=700.1500244140625; E=243.8999938964844; F=474.75; G=1411.849975585938; H=28242.6496887207; I=-0.01601098691069303; J=252.4755664530526 },
    @{ Row=622; Date="2024-09-10"; C=1035.800048828125; D=713.4000244140625; E=248.25; F=478.7999877929688; G=1424.449951171875; H=28522.85046386719; I=0.009921192885042528; J=254.9804252465937 },
    @{ Row=623; Date="2024-09-11"; C=976.2999877929688; D=725.4000244140625; E=241.5500030517578; F=472.2000122070312; G=1399.599975585938; H=27922.25028991699; I=-0.02105680758348599; J=249.6113514946207 },
    @{ Row=624; Date="2024-09-12"; C=986.1500244140625; D=726.0499877929688; E=246.1499938964844; F=479.8500061035156; G=1403.150024414062; H=28182.30001831055; I=0.009313351384414074; J=251.9360697206286 },
    @{ Row=625; Date="2024-09-13"; C=992.0999755859375; D=724.25; E=245.6499938964844; F=485.3999938964844; G=1410.949951171875; H=28258.59951782227; I=0.002707355306775728; J=252.618150175955 },
    @{ Row=626; Date="2024-09-16"; C=988.4000244140625; D=733.6500244140625; E=243.8000030517578; F=489.9500122070312; G=1404.550048828125; H=28313.45072937012; I=0.00194104493795801; J=253.1084933575904 },
    @{ Row=627; Date="2024-09-17"; C=974.9500122070312; D=745.4000244140625; E=240.8000030517578; F=482.2999877929688; G=1400.25; H=28196.30033874512; I=-0.004137623200533361; J=252.0612257832219 },
    @{ Row=628; Date="2024-09-18"; C=962.0499877929688; D=717.5499877929688; E=235.9499969482422; F=471.75; G=1391.300048828125; H=27572.89979553223; I=-0.02210930284198537; J=246.4883278076586 },
    @{ Row=629; Date="2024-09-19"; C=967; D=728.5; E=237.5500030517578; F=459.9500122070312; G=1374.150024414062; H=27641.0502166748; I=0.002471645044516532; J=247.0975594616156 },
    @{ Row=630; Date="2024-09-20"; C=970.8499755859375; D=748.3499755859375; E=237.8500061035156; F=466.2999877929688; G=1380.550048828125; H=27960.69967651367; I=0.01156430227264067; J=249.9550703300615 },
    @{ Row=631; Date="2024-09-23"; C=971.7999877929688; D=750.2000122070312; E=236.4499969482422; F=471.1499938964844; G=1375.400024414062; H=27984.94996643066; I=0.0008672991090191444; J=250.1718561398536 },
    @{ Row=632; Date="2024-09-24"; C=977.2999877929688; D=735.9000244140625; E=237.3000030517578; F=476.7000122070312; G=1363.699951171875; H=27912.20024108887; I=-0.002599601765558408; J=249.5215089409394 },
    @{ Row=633; Date="2024-09-25"; C=963.5999755859375; D=730.0499877929688; E=238.3500061035156; F=473.7000122070312; G=1365.400024414062; H=27741.79995727539; I=-0.006104867489544392; J=247.9982131930638 }
)

foreach ($row in $newRows) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value = $row.Date
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
}
